$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column R (Disponible) from "SI" to "NO" for rows 26 through 47
# because the unit distribution by rounds is now finished for these rows.
for ($row = 26; $row -le 47; $row++) {
    $ws.Range("R$row").Value = "NO"
}

# Update the view state of the sheet to reflect the new working position
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("T27").Select()
